# Generate Report for Handback
# - Row 3 (the e5d2d437-... file) status flips from "Ready for handoff" to
#   "Handback transform failed" on the Overview sheet as well as on each
#   per-language sheet (the same shared text appears in several columns).
# - Each per-language sheet gets a new "Error Detail" (column K) value on
#   row 3 describing the handback/handoff file-name mismatch.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Status column cells for the e5d2d437-...md row (row 3) on every sheet that
# shows it.
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# New "Error Detail" entries (column K) for row 3 on the language sheets.
$zhcn.Range("K3").Value = "Handback file name: vr4xvtie.qf5 is different with handoff file name: e5d2d437-556f-40bb-a3e3-c58b1c66ff90.3e78a54fda66066ca61e18fc4cc00853fd5ed81a.zh-cn."
$dede.Range("K3").Value = "Handback file name: vr4xvtie.qf5 is different with handoff file name: e5d2d437-556f-40bb-a3e3-c58b1c66ff90.3e78a54fda66066ca61e18fc4cc00853fd5ed81a.de-de."
